# "version add product - tax"
# Adds a new "vat" column (M) to the AddProduct sheet with values for each
# product row, and updates the sheet view's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in M1, using the next shared-string slot ("vat").
$ws.Range("M1").Value = "vat"

# New vat values for the three product rows.
$ws.Range("M2").Value = 5
$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 2

# Update the sheet view: scroll so column B is left-most and select N2,
# matching the author's recorded view state after adding the column.
$null = $ws.Range("N2").Select()
$excel.ActiveWindow.ScrollColumn = 2
